$d = $word.ActiveDocument

# 1. "For the last section write up" -> "For the last section, write up"
$d.Content.Find.Execute("For the last section write up", $true, $false, $false, $false, $false,
                         $true, 1, $false, "For the last section, write up", 2)

# 2. Fix redundant/garbled sentence about choosing sections
$d.Content.Find.Execute(
  "we all chose different sections so there, meaning we all wanted to work on different sections so there were no disagreements in the group.",
  $true, $false, $false, $false, $false,
  $true, 1, $false,
  "we all chose different sections, meaning there were no disagreements in the group.", 2)

# 3. "a files tab and notes tab" -> "a “Files” tab and a “Notes” tab"
$d.Content.Find.Execute("a files tab and notes tab", $true, $false, $false, $false, $false,
                         $true, 1, $false, "a “Files” tab and a “Notes” tab", 2)

# 4. "we decided have all members in a group call" -> "we decided to have all members in a scheduled group call"
$d.Content.Find.Execute("we decided have all members in a group call", $true, $false, $false, $false, $false,
                         $true, 1, $false, "we decided to have all members in a scheduled group call", 2)

# 5. "we unmute out microphone" -> "we unmute our microphone"
$d.Content.Find.Execute("we unmute out microphone", $true, $false, $false, $false, $false,
                         $true, 1, $false, "we unmute our microphone", 2)

# 6. "speak to the group that way we get" -> "speak to the group, that way we get"
$d.Content.Find.Execute("speak to the group that way we get", $true, $false, $false, $false, $false,
                         $true, 1, $false, "speak to the group, that way we get", 2)

# 7. "update the repository we do by pushing" -> "update the repository, we do by pushing"
$d.Content.Find.Execute("update the repository we do by pushing", $true, $false, $false, $false, $false,
                         $true, 1, $false, "update the repository, we do by pushing", 2)
